# Atualização automática de TUCUNDUVA.xlsx
$wb = $excel.ActiveWorkbook

# 1. Delete the "Desarquivamentos Pendentes" sheet entirely.
$excel.DisplayAlerts = $false
$ws7 = $wb.Worksheets.Item("Desarquivamentos Pendentes")
$ws7.Delete() | Out-Null

# 2. Rename "Paineis DARQ" -> "PAINEIS DARQ"
$ws1 = $wb.Worksheets.Item("Paineis DARQ")
$ws1.Name = "PAINEIS DARQ"

# 3. Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
$ws5 = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$ws5.Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Keep the original active/selected sheet ("Paineis DARQ" -> "PAINEIS DARQ")
$ws1.Activate() | Out-Null
$ws1.Range("S7").Select() | Out-Null
